# Replace the two-paragraph "{% for citation in literature %}" / "{{citation}}{% endfor %}"
# block in the Bibliography section with the new set of explicit literature.* reference
# paragraphs (citation1 style), as described in the commit "added references to templates".

$d = $word.ActiveDocument

$startMarker = "{% for citation in literature %}"
$endMarker = "{% endfor %}"

# Locate the start paragraph ("{% for citation in literature %}") and the end paragraph
# (the one containing "{{citation}}{% endfor %}") by scanning the document paragraphs.
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains($startMarker)) {
        $startIndex = $i
    }
    if ($startIndex -ge 0 -and $t.Contains($endMarker)) {
        $endIndex = $i
        break
    }
}

if ($startIndex -lt 0 -or $endIndex -lt 0) {
    throw "Could not locate the citation for/endfor paragraphs"
}

$rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
$rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
$target = $d.Range($rangeStart, $rangeEnd)

$newParagraphsXml = '<w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.integration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.absorption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.</w:t></w:r><w:r><w:t>solution</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.refinement</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.</w:t></w:r><w:r><w:t>ccdc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>literature.</w:t></w:r><w:r><w:t>finalcif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="360"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="citation1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:left="360"/></w:pPr></w:p>'

$wordXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document ' + $wordXmlNs + '><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData>' + `
    '</pkg:part></pkg:package>'

$paragraphCountBefore = $d.Paragraphs.Count

$target.InsertXML($pkgXml)

$paragraphCountAfter = $d.Paragraphs.Count

# InsertXML (via the underlying list-style reconciliation logic) silently drops an
# explicit <w:ind w:left="360"/> that coincides with the list level's own indent value
# when paired with a <w:numPr> that cancels numbering (numId=0). Re-apply the indent
# via ParagraphFormat on the two trailing empty paragraphs so the OOXML round-trips
# with the explicit w:ind element, matching the target markup.
$lastParaIndex = $paragraphCountAfter
$secondLastParaIndex = $paragraphCountAfter - 1
$d.Paragraphs.Item($secondLastParaIndex).Range.ParagraphFormat.LeftIndent = 18
$d.Paragraphs.Item($lastParaIndex).Range.ParagraphFormat.LeftIndent = 18

Write-Host "Replaced citation loop block with explicit literature references"
Write-Host "Paragraphs before: $paragraphCountBefore, after: $paragraphCountAfter"
